$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45206 -> 2023-10-07).
# Update it to serial 45208 (2023-10-09) for all data rows (2 through 519).
$ws.Range("C2:C519").Value = 45208
